$d = $word.ActiveDocument

# --- Change 1 ---
# "...thời gian hiện tại hay không." -> "...thời gian hiện tại hay không và kiểm tra số lượng nhóm còn có thể đăng ký."
$d.Content.Find.Execute(
    "Hệ thống sẽ kiểm tra thời gian nộp của đồ án đó xem có lớn hơn thời gian hiện tại hay không.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Hệ thống sẽ kiểm tra thời gian nộp của đồ án đó xem có lớn hơn thời gian hiện tại hay không và kiểm tra số lượng nhóm còn có thể đăng ký.",
    2) | Out-Null

# --- Change 2 ---
# "thời gian nộp lớn hơn thời gian hiện tại " -> "thỏa các điều kiện nói trên "
$d.Content.Find.Execute(
    "thời gian nộp lớn hơn thời gian hiện tại ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "thỏa các điều kiện nói trên ",
    2) | Out-Null

# --- Change 3 ---
# "...của một đồ án cụ thể." -> "...của một đồ án thuộc một môn học."
$d.Content.Find.Execute(
    "Giảng viên có nhu cầu cập nhật lại thời gian nộp của một đồ án cụ thể.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Giảng viên có nhu cầu cập nhật lại thời gian nộp của một đồ án thuộc một môn học.",
    2) | Out-Null

# --- Change 4 (part 1) ---
# The paragraph's text is replaced by the new "kiểm tra giáo viên..." sentence.
$d.Content.Find.Execute(
    "Hệ thống sẽ cập nhật thời gian của đồ án trong bảng DOAN với thời gian nộp muốn chỉnh sửa và mã đồ án.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Hệ thống sẽ kiểm tra giáo viên có phụ trách môn học này không và thời gian nộp mới có trước thời gian kết thúc môn học hay không.",
    2) | Out-Null

# --- Change 4 (part 2) ---
# A brand-new bullet paragraph is appended after the paragraph above, carrying the
# sentence that used to live there before the edit ("Hệ thống sẽ cập nhật ... mã đồ án."),
# now prefixed with "Nếu thỏa các điều kiện nói trên thì ".
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "Nếu thỏa các điều kiện nói trên thì hệ thống sẽ cập nhật thời gian của đồ án trong bảng DOAN với thời gian nộp muốn chỉnh sửa và mã đồ án."
